$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-05 18:55:42"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
